$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3454.9119
$ws.Range("I15").Value = 3454.9119
$ws.Range("K15").Value = 10364.7357
$ws.Range("M15").Value = -10195.7357
$ws.Range("H18").Value = 710
$ws.Range("I18").Value = 710
$ws.Range("K18").Value = 710
$ws.Range("M18").Value = -426
$ws.Range("H100").Value = 5159.8
$ws.Range("I100").Value = 4933
$ws.Range("K100").Value = 4933
$ws.Range("M100").Value = -4392
$ws.Range("H132").Value = 10902.523
$ws.Range("I132").Value = 10902.523
$ws.Range("K132").Value = 32707.569
$ws.Range("M132").Value = -30177.569
$ws.Range("H137").Value = 2705.5334
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 2755.9285
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 8267.7855
$ws.Range("M137").Value = -3450
$ws.Range("N137").Value = -13367.7855
$ws.Range("H138").Value = 8033.143
$ws.Range("J138").Value = 9288.666999999999
$ws.Range("L138").Value = 27866.001
$ws.Range("N138").Value = -38146.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2199.6
$ws.Range("I2").Value = 2199.6
$ws.Range("K2").Value = 2199.6
$ws.Range("M2").Value = -2086.6
$ws.Range("H17").Value = 2008.3334
$ws.Range("I17").Value = 2500
$ws.Range("J17").Value = 1516.6666
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 1516.6666
$ws.Range("M17").Value = -2327
$ws.Range("N17").Value = -1862.6666
$ws.Range("H74").Value = 1706.8846
$ws.Range("I74").Value = 1465.5416
$ws.Range("J74").Value = 4603
$ws.Range("K74").Value = 1465.5416
$ws.Range("L74").Value = 4603
$ws.Range("M74").Value = -591.5416
$ws.Range("N74").Value = -6351
$ws.Range("H77").Value = 1706.8846
$ws.Range("I77").Value = 1465.5416
$ws.Range("J77").Value = 4603
$ws.Range("K77").Value = 7327.708000000001
$ws.Range("L77").Value = 23015
$ws.Range("M77").Value = -2959.708000000001
$ws.Range("N77").Value = -31751
$ws.Range("H88").Value = 1727.4445
$ws.Range("I88").Value = 1563.5
$ws.Range("K88").Value = 1563.5
$ws.Range("M88").Value = -1157.5
$ws.Range("H91").Value = 1727.4445
$ws.Range("I91").Value = 1563.5
$ws.Range("K91").Value = 1563.5
$ws.Range("M91").Value = -159.5
$ws.Range("H97").Value = 1008.1539
$ws.Range("I97").Value = 945.9091
$ws.Range("K97").Value = 945.9091
$ws.Range("M97").Value = -449.9091
$ws.Range("H110").Value = 166670780
$ws.Range("I110").Value = 250003300
$ws.Range("K110").Value = 250003300
$ws.Range("M110").Value = -250001255
$ws.Range("H116").Value = 2199.6
$ws.Range("I116").Value = 2199.6
$ws.Range("K116").Value = 2199.6
$ws.Range("M116").Value = 94.40000000000009
$ws.Range("H132").Value = 1446.5238
$ws.Range("I132").Value = 1398.7222
$ws.Range("K132").Value = 4196.1666
$ws.Range("M132").Value = -1666.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2199.6
$ws.Range("I3").Value = 2199.6
$ws.Range("K3").Value = 2199.6
$ws.Range("M3").Value = -2085.6
$ws.Range("H105").Value = 7396609.5
$ws.Range("I105").Value = 12326016
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 12326016
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = -12324269
$ws.Range("N105").Value = -5994

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2033
$ws.Range("I16").Value = 550
$ws.Range("J16").Value = 4999
$ws.Range("K16").Value = 550
$ws.Range("L16").Value = 4999
$ws.Range("M16").Value = -263
$ws.Range("N16").Value = -5573
$ws.Range("H31").Value = 5692.2334
$ws.Range("I31").Value = 2848.1667
$ws.Range("K31").Value = 2848.1667
$ws.Range("M31").Value = -2553.1667
$ws.Range("H34").Value = 5692.2334
$ws.Range("I34").Value = 2848.1667
$ws.Range("K34").Value = 2848.1667
$ws.Range("M34").Value = -2646.1667
$ws.Range("H58").Value = 1392.5
$ws.Range("I58").Value = 732.2222
$ws.Range("K58").Value = 732.2222
$ws.Range("M58").Value = -529.2222
$ws.Range("H68").Value = 86228.336
$ws.Range("J68").Value = 86228.336
$ws.Range("L68").Value = 86228.336
$ws.Range("N68").Value = -87726.336
$ws.Range("H71").Value = 86228.336
$ws.Range("J71").Value = 86228.336
$ws.Range("L71").Value = 258685.008
$ws.Range("N71").Value = -266173.008
$ws.Range("H113").Value = 2033
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 4999
$ws.Range("K113").Value = 550
$ws.Range("L113").Value = 4999
$ws.Range("M113").Value = 1620
$ws.Range("N113").Value = -9339
$ws.Range("H132").Value = 1264.5927
$ws.Range("J132").Value = 24
$ws.Range("L132").Value = 72
$ws.Range("N132").Value = -5132
$ws.Range("H136").Value = 1392.5
$ws.Range("I136").Value = 732.2222
$ws.Range("K136").Value = 2196.6666
$ws.Range("M136").Value = 353.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 296.75
$ws.Range("I50").Value = 435
$ws.Range("J50").Value = 66.333336
$ws.Range("K50").Value = 1305
$ws.Range("L50").Value = 199.000008
$ws.Range("M50").Value = -824
$ws.Range("N50").Value = -1161.000008
$ws.Range("H53").Value = 296.75
$ws.Range("I53").Value = 435
$ws.Range("J53").Value = 66.333336
$ws.Range("K53").Value = 1305
$ws.Range("L53").Value = 199.000008
$ws.Range("M53").Value = -824
$ws.Range("N53").Value = -1161.000008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H62").Value = 2500000
$ws.Range("I62").Value = 2500000
$ws.Range("K62").Value = 2500000
$ws.Range("M62").Value = -2499314
$ws.Range("H65").Value = 2500000
$ws.Range("I65").Value = 2500000
$ws.Range("K65").Value = 7500000
$ws.Range("M65").Value = -7496568
$ws.Range("H122").Value = 3498.5715
$ws.Range("J122").Value = 3653.25
$ws.Range("L122").Value = 10959.75
$ws.Range("N122").Value = -15859.75
$ws.Range("H132").Value = 40442.652
$ws.Range("I132").Value = 43562.875
$ws.Range("K132").Value = 130688.625
$ws.Range("M132").Value = -128158.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H40").Value = 4122.4346
$ws.Range("I40").Value = 3946.182
$ws.Range("K40").Value = 3946.182
$ws.Range("M40").Value = -3810.182
$ws.Range("H46").Value = 7363.636
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 7363.636
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 7363.636
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -7739.636
$ws.Range("H93").Value = 1558.238
$ws.Range("I93").Value = 1453.1111
$ws.Range("K93").Value = 1453.1111
$ws.Range("M93").Value = -205.1111000000001
$ws.Range("H122").Value = 3134.9092
$ws.Range("I122").Value = 3164.889
$ws.Range("K122").Value = 9494.667000000001
$ws.Range("M122").Value = -7044.667000000001
$ws.Range("H136").Value = 1159.5
$ws.Range("I136").Value = 1159.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3478.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -928.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 37483.332
$ws.Range("J103").Value = 37483.332
$ws.Range("L103").Value = 37483.332
$ws.Range("N103").Value = -39827.332
$ws.Range("H126").Value = 3688.3333
$ws.Range("J126").Value = 8000
$ws.Range("L126").Value = 24000
$ws.Range("N126").Value = -28940
$ws.Range("H136").Value = 3190.7693
$ws.Range("I136").Value = 2110.375
$ws.Range("J136").Value = 4919.4
$ws.Range("K136").Value = 6331.125
$ws.Range("L136").Value = 14758.2
$ws.Range("M136").Value = -3781.125
$ws.Range("N136").Value = -19858.2
